$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (GitHub Actions cron).
# Columns D (Price) and E (Volume(1h)) hold plain-text, pre-formatted
# strings (e.g. "67.30", "  +1.30%  "). Force NumberFormat to Text on
# each touched cell first so Excel does not coerce numeric-looking
# strings into real numbers (which would drop trailing zeros / exponents).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.745.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.441.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.71'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.25%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.94'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.33%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.882.77'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.621.69'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.444.29'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.29'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.93'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.64'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.30'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.82'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.71'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +9.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '576.80'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.561.85'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.05%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.85'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.71%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.43'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '148.05'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.13%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +7.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '148.28'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.67'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.52'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.31%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.53%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0231'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0924'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.41%  '
